$d = $word.ActiveDocument

# 1. Merge "5" + "/2020" into "5/2020" in the version/date cell (table 1, row 2, col 1).
$d.Content.Find.Execute("1/25/2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1/25/2020", 0) | Out-Null

# 2. Merge "1_" + "2" into "1_2" in the REQ id cell.
$d.Content.Find.Execute("_CYRS_06_V1_2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "_CYRS_06_V1_2", 0) | Out-Null

# 3. Change "Logging sequence" -> "Logging system" with "system" bold.
$range = $d.Content
$found = $range.Find.Execute("Logging sequence", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
if ($found) {
    $range.Text = "Logging system"
}
